$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.567946
$ws.Range("H2").Value = 7.703838
$ws.Range("I2").Value = 0.8361514603436194
$ws.Range("J2").Value = 0.8361514603436195
$ws.Range("Q2").Value = 0.160019843026
$ws.Range("R2").Value = 1.440178587234
$ws.Range("S2").Value = 0.8361514603436194
$ws.Range("T2").Value = 0.8361514603436195

# Row 3 updates
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5032033333333333
$ws.Range("H3").Value = 1.50961
$ws.Range("I3").Value = 0.1638485396563805
$ws.Range("J3").Value = 0.1638485396563805
$ws.Range("Q3").Value = 0.03135678024777778
$ws.Range("R3").Value = 0.28221102223
$ws.Range("S3").Value = 0.1638485396563805
$ws.Range("T3").Value = 0.1638485396563805
